# repull data, push all data, mean calculation
# Update the dSF column (F) for several rows with refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -9
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = 2
$ws.Range("F9").Value = -5
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = -8
$ws.Range("F17").Value = -3
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 4
$ws.Range("F23").Value = -2
